$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is inserted at row 351, pushing the existing
# rows 351:368 down to 352:369 (dimension grows from R368 to R369).
$ws.Rows(351).Insert()

# Populate the newly inserted row 351 with the new record.
$ws.Cells.Item(351, 1).Value = 9
$ws.Cells.Item(351, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(351, 3).Value = "Metropolitana"
$ws.Cells.Item(351, 4).Value = 44939
$ws.Cells.Item(351, 5).Value = 13
$ws.Cells.Item(351, 6).Value = 300000001
$ws.Cells.Item(351, 7).Value = "Rabanito"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "Primera"
$ws.Cells.Item(351, 10).Value = 7000
$ws.Cells.Item(351, 11).Value = 3000
$ws.Cells.Item(351, 12).Value = 3000
$ws.Cells.Item(351, 13).Value = 3000
$ws.Cells.Item(351, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(351, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(351, 16).Value = 30
$ws.Cells.Item(351, 17).Value = 100
$ws.Cells.Item(351, 18).Value = "Hortaliza"
